# Apply updated crypto price/volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.367.36"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "2.519.36"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.61"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.69"
$ws.Range("E6").Value = "  -3.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("E8").Value = "  -1.30%  "
$ws.Range("D9").Value = "2.524.31"
$ws.Range("E9").Value = "  -1.72%  "
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.162"
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.43"
$ws.Range("E12").Value = "  -2.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.350"
$ws.Range("E13").Value = "  -3.42%  "
$ws.Range("D14").Value = "2.974.42"
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.28"
$ws.Range("E15").Value = "  -2.53%  "
$ws.Range("D16").Value = "59.296.65"
$ws.Range("E16").Value = "  -0.94%  "
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D18").Value = "2.516.07"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.08"
$ws.Range("E19").Value = "  -2.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.27"
$ws.Range("E20").Value = "  -1.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.98"
$ws.Range("E21").Value = "  -0.92%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("E23").Value = "  -0.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.29"
$ws.Range("E24").Value = "  +1.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.421"
$ws.Range("E25").Value = "  -4.54%  "
$ws.Range("E26").Value = "  +0.72%  "
$ws.Range("E27").Value = "  +0.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.71"
$ws.Range("E28").Value = "  -3.79%  "
$ws.Range("D29").Value = "0.0₃0782"
$ws.Range("E29").Value = "  -2.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.70"
$ws.Range("E30").Value = "  -5.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.79"
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.54"
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.11"
$ws.Range("E34").Value = "  -9.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.40"
$ws.Range("E35").Value = "  -6.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.55"
$ws.Range("E36").Value = "  -1.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.13"
$ws.Range("E37").Value = "  -7.53%  "
$ws.Range("E38").Value = "  -2.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.67"
$ws.Range("E39").Value = "  -1.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.808"
$ws.Range("E40").Value = "  -3.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.22"
$ws.Range("E41").Value = "  -8.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "281.15"
$ws.Range("E42").Value = "  -6.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("E44").Value = "  -1.47%  "
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "126.09"
$ws.Range("E46").Value = "  +1.12%  "
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0512"
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0223"
$ws.Range("E49").Value = "  -2.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.81"
$ws.Range("E50").Value = "  -2.77%  "
$ws.Range("D51").Value = "1.770.54"
$ws.Range("E51").Value = "  -2.76%  "
